$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.843.80"
$ws.Range("E2").Value = "  +1.99%  "
$ws.Range("D3").Value = "2.118.83"
$ws.Range("E3").Value = "  +6.19%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'334.32"
$ws.Range("E5").Value = "  +2.94%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.5328"
$ws.Range("E7").Value = "  +4.56%  "
$ws.Range("D8").Value = "'0.4384"
$ws.Range("E8").Value = "  +6.26%  "
$ws.Range("D9").Value = "'0.09025"
$ws.Range("E9").Value = "  +3.81%  "
$ws.Range("D10").Value = "'46.97"
$ws.Range("E10").Value = "  +9.59%  "
$ws.Range("E11").Value = "  +4.19%  "
$ws.Range("D12").Value = "'25.02"
$ws.Range("E12").Value = "  +2.46%  "
$ws.Range("D13").Value = "2.116.16"
$ws.Range("E13").Value = "  +6.21%  "
$ws.Range("D14").Value = "'6.761"
$ws.Range("E14").Value = "  +4.07%  "
$ws.Range("D15").Value = "'7.805"
$ws.Range("E15").Value = "  +5.30%  "
$ws.Range("D16").Value = "'97.17"
$ws.Range("E16").Value = "  +3.46%  "
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "'0.00001138"
$ws.Range("E18").Value = "  +2.23%  "
$ws.Range("D19").Value = "'0.06679"
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("D20").Value = "'19.13"
$ws.Range("E20").Value = "  +1.94%  "
$ws.Range("D21").Value = "'0.9999"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "'6.344"
$ws.Range("E22").Value = "  +4.19%  "
$ws.Range("D23").Value = "30.914.63"
$ws.Range("E23").Value = "  +2.06%  "
$ws.Range("D24").Value = "'12.27"
$ws.Range("E24").Value = "  +4.83%  "
$ws.Range("D25").Value = "2.365.46"
$ws.Range("E25").Value = "  +6.35%  "
$ws.Range("D26").Value = "'2.293"
$ws.Range("E26").Value = "  +3.71%  "
$ws.Range("D27").Value = "'22.75"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("D28").Value = "'2.617"
$ws.Range("E28").Value = "  +10.45%  "
$ws.Range("D29").Value = "'163.41"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "'133.54"
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("D31").Value = "'1.193"
$ws.Range("E31").Value = "  +5.67%  "
$ws.Range("E32").Value = "  +3.28%  "
$ws.Range("D33").Value = "'6.221"
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("D34").Value = "'4.053"
$ws.Range("E34").Value = "  +6.28%  "
$ws.Range("D35").Value = "'1.556"
$ws.Range("E35").Value = "  +18.24%  "
$ws.Range("D36").Value = "'0.02603"
$ws.Range("E36").Value = "  +5.12%  "
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "'9.630"
$ws.Range("E37").Value = "  +8.36%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06793"
$ws.Range("E38").Value = "  +3.86%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.549"
$ws.Range("E39").Value = "  +3.13%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2321"
$ws.Range("E40").Value = "  +5.84%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'12.74"
$ws.Range("E41").Value = "  +8.51%  "
$ws.Range("D42").Value = "'0.6862"
$ws.Range("E42").Value = "  +4.26%  "
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("D44").Value = "'0.6478"
$ws.Range("E44").Value = "  +5.90%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'14.07"
$ws.Range("E46").Value = "  +3.65%  "
$ws.Range("D47").Value = "'2.241"
$ws.Range("E47").Value = "  +2.20%  "
$ws.Range("D48").Value = "'3.665"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").Value = "'1.269"
$ws.Range("E49").Value = "  +3.71%  "
$ws.Range("B50").Value = "WEMIXTOKEN"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").Value = "'1.205"
$ws.Range("E50").Value = "  +11.17%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'83.08"
$ws.Range("E51").Value = "  +4.65%  "
